$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 134
$ws.Range("I4").Value = 75
$ws.Range("J4").Value = 173.33333
$ws.Range("K4").Value = 75
$ws.Range("L4").Value = 173.33333
$ws.Range("M4").Value = 39
$ws.Range("N4").Value = -401.33333
$ws.Range("H9").Value = 8422.462
$ws.Range("I9").Value = 11937.556
$ws.Range("J9").Value = 513.5
$ws.Range("K9").Value = 11937.556
$ws.Range("L9").Value = 513.5
$ws.Range("M9").Value = -11768.556
$ws.Range("N9").Value = -851.5
$ws.Range("H17").Value = 696880.6
$ws.Range("I17").Value = 1852.3334
$ws.Range("J17").Value = 827198.4399999999
$ws.Range("K17").Value = 5557.0002
$ws.Range("L17").Value = 2481595.32
$ws.Range("M17").Value = -5389.0002
$ws.Range("N17").Value = -2481931.32
$ws.Range("H64").Value = 4666.6665
$ws.Range("J64").Value = 5500
$ws.Range("L64").Value = 5500
$ws.Range("N64").Value = -5996
$ws.Range("H67").Value = 4666.6665
$ws.Range("J67").Value = 5500
$ws.Range("L67").Value = 5500
$ws.Range("N67").Value = -7216
$ws.Range("H107").Value = 701.9091
$ws.Range("J107").Value = 491.5
$ws.Range("L107").Value = 491.5
$ws.Range("N107").Value = -4331.5
$ws.Range("H132").Value = 3718.9443
$ws.Range("I132").Value = 3558.875
$ws.Range("K132").Value = 10676.625
$ws.Range("M132").Value = -8146.625
$ws.Range("H141").Value = 2128.1333
$ws.Range("I141").Value = 2302.077
$ws.Range("J141").Value = 997.5
$ws.Range("K141").Value = 6906.231000000001
$ws.Range("L141").Value = 2992.5
$ws.Range("M141").Value = -1726.231000000001
$ws.Range("N141").Value = -13352.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1851.6904
$ws.Range("I32").Value = 1861.775
$ws.Range("J32").Value = 1650
$ws.Range("K32").Value = 1861.775
$ws.Range("L32").Value = 1650
$ws.Range("M32").Value = -1574.775
$ws.Range("N32").Value = -2224
$ws.Range("H45").Value = 5708.3
$ws.Range("I45").Value = 5583.7144
$ws.Range("K45").Value = 5583.7144
$ws.Range("M45").Value = -5206.7144
$ws.Range("H74").Value = 20835766
$ws.Range("I74").Value = 25001944
$ws.Range("K74").Value = 25001944
$ws.Range("M74").Value = -25001070
$ws.Range("H77").Value = 20835766
$ws.Range("I77").Value = 25001944
$ws.Range("K77").Value = 125009720
$ws.Range("M77").Value = -125005352
$ws.Range("H106").Value = 65370
$ws.Range("J106").Value = 65370
$ws.Range("L106").Value = 65370
$ws.Range("N106").Value = -67894
$ws.Range("H110").Value = 202879.6
$ws.Range("I110").Value = 334799.34
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 334799.34
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = -332754.34
$ws.Range("N110").Value = -9090
$ws.Range("H122").Value = 5740.4517
$ws.Range("I122").Value = 3652.1155
$ws.Range("J122").Value = 16599.8
$ws.Range("K122").Value = 10956.3465
$ws.Range("L122").Value = 49799.39999999999
$ws.Range("M122").Value = -8506.3465
$ws.Range("N122").Value = -54699.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 632.2
$ws.Range("I5").Value = 739
$ws.Range("J5").Value = 205
$ws.Range("K5").Value = 739
$ws.Range("L5").Value = 205
$ws.Range("M5").Value = -626
$ws.Range("N5").Value = -431
$ws.Range("H20").Value = 651.3333
$ws.Range("I20").Value = 601.8
$ws.Range("J20").Value = 899
$ws.Range("K20").Value = 601.8
$ws.Range("L20").Value = 899
$ws.Range("M20").Value = -354.8
$ws.Range("N20").Value = -1393
$ws.Range("H107").Value = 202921.6
$ws.Range("I107").Value = 1536.3334
$ws.Range("J107").Value = 504999.5
$ws.Range("K107").Value = 1536.3334
$ws.Range("L107").Value = 504999.5
$ws.Range("M107").Value = 383.6666
$ws.Range("N107").Value = -508839.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 7952.615
$ws.Range("I7").Value = 11441.111
$ws.Range("J7").Value = 103.5
$ws.Range("K7").Value = 11441.111
$ws.Range("L7").Value = 103.5
$ws.Range("M7").Value = -11328.111
$ws.Range("N7").Value = -329.5
$ws.Range("H16").Value = 1360980
$ws.Range("I16").Value = 1813585.6
$ws.Range("K16").Value = 1813585.6
$ws.Range("M16").Value = -1813298.6
$ws.Range("H31").Value = 5830.467
$ws.Range("I31").Value = 3573.6155
$ws.Range("K31").Value = 3573.6155
$ws.Range("M31").Value = -3278.6155
$ws.Range("H34").Value = 5830.467
$ws.Range("I34").Value = 3573.6155
$ws.Range("K34").Value = 3573.6155
$ws.Range("M34").Value = -3371.6155
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H58").Value = 19239646
$ws.Range("I58").Value = 20009184
$ws.Range("K58").Value = 20009184
$ws.Range("M58").Value = -20008981
$ws.Range("H60").Value = 10538.308
$ws.Range("J60").Value = 31499.666
$ws.Range("L60").Value = 31499.666
$ws.Range("N60").Value = -32521.666
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H99").Value = 3337.25
$ws.Range("I99").Value = 3242.8572
$ws.Range("K99").Value = 3242.8572
$ws.Range("M99").Value = -1744.8572
$ws.Range("H113").Value = 1360980
$ws.Range("I113").Value = 1813585.6
$ws.Range("K113").Value = 1813585.6
$ws.Range("M113").Value = -1811415.6
$ws.Range("H126").Value = 3337.25
$ws.Range("I126").Value = 3242.8572
$ws.Range("K126").Value = 9728.571599999999
$ws.Range("M126").Value = -7258.571599999999
$ws.Range("H134").Value = 12502594
$ws.Range("I134").Value = 13891224
$ws.Range("K134").Value = 41673672
$ws.Range("M134").Value = -41671137
$ws.Range("H136").Value = 19239646
$ws.Range("I136").Value = 20009184
$ws.Range("K136").Value = 60027552
$ws.Range("M136").Value = -60025002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 784.8461
$ws.Range("I23").Value = 100
$ws.Range("J23").Value = 990.3
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 2970.9
$ws.Range("M23").Value = -65
$ws.Range("N23").Value = -3440.9
$ws.Range("H113").Value = 111715
$ws.Range("I113").Value = 333586.34
$ws.Range("J113").Value = 779.3333
$ws.Range("K113").Value = 1000759.02
$ws.Range("L113").Value = 2337.9999
$ws.Range("M113").Value = -998589.02
$ws.Range("N113").Value = -6677.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8249.5
$ws.Range("I70").Value = 8332.666999999999
$ws.Range("K70").Value = 8332.666999999999
$ws.Range("M70").Value = -8062.666999999999
$ws.Range("H73").Value = 8249.5
$ws.Range("I73").Value = 8332.666999999999
$ws.Range("K73").Value = 8332.666999999999
$ws.Range("M73").Value = -7396.666999999999
$ws.Range("H102").Value = 4379.4
$ws.Range("I102").Value = 4379.4
$ws.Range("K102").Value = 4379.4
$ws.Range("M102").Value = -2757.4
$ws.Range("H113").Value = 22944.08
$ws.Range("I113").Value = 29020.77
$ws.Range("J113").Value = 1399.4546
$ws.Range("K113").Value = 29020.77
$ws.Range("L113").Value = 1399.4546
$ws.Range("M113").Value = -26850.77
$ws.Range("N113").Value = -5739.4546
$ws.Range("H122").Value = 5012.95
$ws.Range("I122").Value = 2903.2778
$ws.Range("K122").Value = 8709.8334
$ws.Range("M122").Value = -6259.8334
$ws.Range("H132").Value = 3575765
$ws.Range("I132").Value = 3575765
$ws.Range("K132").Value = 10727295
$ws.Range("M132").Value = -10724765

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2497
$ws.Range("I7").Value = 1999
$ws.Range("J7").Value = 2995
$ws.Range("K7").Value = 1999
$ws.Range("L7").Value = 2995
$ws.Range("M7").Value = -1887
$ws.Range("N7").Value = -3219
$ws.Range("H22").Value = 2924.4783
$ws.Range("I22").Value = 2498.6
$ws.Range("J22").Value = 3252.077
$ws.Range("K22").Value = 2498.6
$ws.Range("L22").Value = 3252.077
$ws.Range("M22").Value = -2203.6
$ws.Range("N22").Value = -3842.077
$ws.Range("H27").Value = 2924.4783
$ws.Range("I27").Value = 2498.6
$ws.Range("J27").Value = 3252.077
$ws.Range("K27").Value = 2498.6
$ws.Range("L27").Value = 3252.077
$ws.Range("M27").Value = -2391.6
$ws.Range("N27").Value = -3466.077
$ws.Range("H47").Value = 34065
$ws.Range("J47").Value = 34065
$ws.Range("L47").Value = 34065
$ws.Range("N47").Value = -35045
$ws.Range("H52").Value = 34065
$ws.Range("J52").Value = 34065
$ws.Range("L52").Value = 34065
$ws.Range("N52").Value = -34531
$ws.Range("H126").Value = 2497
$ws.Range("I126").Value = 1999
$ws.Range("J126").Value = 2995
$ws.Range("K126").Value = 5997
$ws.Range("L126").Value = 8985
$ws.Range("M126").Value = -3527
$ws.Range("N126").Value = -13925
$ws.Range("H132").Value = 15650936
$ws.Range("I132").Value = 20854666
$ws.Range("K132").Value = 62563998
$ws.Range("M132").Value = -62561468
$ws.Range("H136").Value = 1526.8846
$ws.Range("I136").Value = 1388.6
$ws.Range("K136").Value = 4165.799999999999
$ws.Range("M136").Value = -1615.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1658
$ws.Range("I113").Value = 1716.4546
$ws.Range("J113").Value = 1497.25
$ws.Range("K113").Value = 5149.3638
$ws.Range("L113").Value = 4491.75
$ws.Range("M113").Value = -2979.3638
$ws.Range("N113").Value = -8831.75
$ws.Range("H136").Value = 21741618
$ws.Range("I136").Value = 21741618
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 65224854
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -65222304
$ws.Range("N136").ClearContents()
